$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase the "Số ngày" (days) value for the Code module activities (rows 10-16)
# from 1 to 2
$ws.Range("C10:C16").Value = 2

# Update the selected cell to C16 to match the final cursor position
$ws.Range("C16").Select()
